$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# rao_hmec_10kb (sheet4): new accuracy numbers + mtry-tuning block
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("rao_hmec_10kb")

$ws4.Range("C3").Value = 0.6802
$ws4.Range("D3").Value = 0.6718
$ws4.Range("C4").Value = 0.6241
$ws4.Range("D4").Value = 0.6017
$ws4.Range("C5").Value = 0.7363
$ws4.Range("D5").Value = 0.7419
$ws4.Range("C6").Value = 0.6612
$ws4.Range("D6").Value = 0.6471
$ws4.Range("C7").Value = 0.703
$ws4.Range("D7").Value = 0.6998
$ws4.Range("C8").Value = 0.6241
$ws4.Range("D8").Value = 0.6017

$ws4.Range("C12").Value = "total 111 predictors"
$ws4.Range("C13").Value = "mtry "
$ws4.Range("C14").Value = "2:112 by 2"
$ws4.Range("C15").Value = "final value used 74"
$ws4.Range("C16").Value = "accuracy at 74 is 0.6783"
$ws4.Rows.Item(16).RowHeight = 26

$ws4.Activate()
$ws4.Range("A12:C20").Select()

# ---------------------------------------------------------------------------
# rao_hela_10kb (sheet3): new accuracy numbers + mtry-tuning block
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("rao_hela_10kb")

$ws3.Range("C3").Value = 0.7148
$ws3.Range("C4").Value = 0.6739
$ws3.Range("C5").Value = 0.7557
$ws3.Range("C6").Value = 0.7026
$ws3.Range("C7").Value = 0.734
$ws3.Range("C8").Value = 0.6739

$ws3.Range("C13").Value = "total 385 predictors"
$ws3.Range("C14").Value = "mtry "
$ws3.Range("C15").Value = "9:378 by 9"
$ws3.Range("C16").Value = "final value used 252"
$ws3.Range("C17").Value = "accuracy at 252 is 0.697127"
$ws3.Rows.Item(17).RowHeight = 33

$ws3.Activate()
$ws3.Range("C9").Select()

# ---------------------------------------------------------------------------
# rao_huvec_10kb (sheet5): new accuracy numbers + training/testing + mtry block
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("rao_huvec_10kb")

$ws5.Range("C3").Value = 0.712
$ws5.Range("C4").Value = 0.6455
$ws5.Range("C5").Value = 0.7784
$ws5.Range("C6").Value = 0.6915
$ws5.Range("C7").Value = 0.7445
$ws5.Range("C8").Value = 0.6455

$ws5.Range("A18").Value = "training"
$ws5.Range("B18").Value = "5416 obs of 369 var"
$ws5.Range("A19").Value = "testing"
$ws5.Range("B19").Value = "1354 obs of 369 var"

$ws5.Range("C12").Value = "total 368 predictors"
$ws5.Range("C13").Value = "mtry "
$ws5.Range("C14").Value = "9:369 by 9"
$ws5.Range("C15").Value = "final value used 297"
$ws5.Range("C16").Value = "accuracy at 74 is 0.71288"
$ws5.Rows.Item(16).RowHeight = 31

$ws5.Activate()
$ws5.Range("C9").Select()

# ---------------------------------------------------------------------------
# dixon_h1esc  (sheet2): new "positive class" / "random" note, becomes active tab
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("dixon_h1esc ")

$ws2.Range("A23").Value = "positive class"
$ws2.Range("B23").Value = "random"

$ws2.Activate()
$ws2.Range("E19").Select()
